$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 0.6077856373797204
$ws.Cells.Item(2, 4).Value2 = 0.08858407273058155
$ws.Cells.Item(2, 5).Value2 = 0.1383886882530991
$ws.Cells.Item(2, 6).Value2 = 2.425425980498446
$ws.Cells.Item(2, 7).Value2 = 1.739103674548488
$ws.Cells.Item(2, 8).Value2 = 1.533532072520757
$ws.Cells.Item(2, 9).Value2 = 1.661324305719521
$ws.Cells.Item(2, 10).Value2 = 0.2178276543611997
$ws.Cells.Item(2, 11).Value2 = 0.501783247074485
$ws.Cells.Item(2, 12).Value2 = 0.2826155600626947
$ws.Cells.Item(2, 13).Value2 = 0.2062002548140995
$ws.Cells.Item(3, 2).Value2 = 0.5905500464506019
$ws.Cells.Item(3, 4).Value2 = 0.08852275330679227
$ws.Cells.Item(3, 5).Value2 = 0.1385762266425049
$ws.Cells.Item(3, 6).Value2 = 2.429558851590372
$ws.Cells.Item(3, 7).Value2 = 1.739172392453639
$ws.Cells.Item(3, 8).Value2 = 1.538815280680907
$ws.Cells.Item(3, 9).Value2 = 1.671441036884428
$ws.Cells.Item(3, 10).Value2 = 0.2179082631518376
$ws.Cells.Item(3, 11).Value2 = 0.4535836331517089
$ws.Cells.Item(3, 12).Value2 = 0.2759665249228505
$ws.Cells.Item(3, 13).Value2 = 0.2007994489899332
$ws.Cells.Item(4, 2).Value2 = 0.5802164792265927
$ws.Cells.Item(4, 4).Value2 = 0.08849622883738384
$ws.Cells.Item(4, 5).Value2 = 0.1387009929940231
$ws.Cells.Item(4, 6).Value2 = 2.433141834604491
$ws.Cells.Item(4, 7).Value2 = 1.740084023668345
$ws.Cells.Item(4, 8).Value2 = 1.542651244385269
$ws.Cells.Item(4, 9).Value2 = 1.678229608420406
$ws.Cells.Item(4, 10).Value2 = 0.2179665382503924
$ws.Cells.Item(4, 11).Value2 = 0.424123304830573
$ws.Cells.Item(4, 12).Value2 = 0.2720137689323252
$ws.Cells.Item(4, 13).Value2 = 0.1975723632342508
$ws.Cells.Item(5, 2).Value2 = 0.576068547328461
$ws.Cells.Item(5, 4).Value2 = 0.08848822963513214
$ws.Cells.Item(5, 5).Value2 = 0.1387542589701244
$ws.Cells.Item(5, 6).Value2 = 2.434864944551585
$ws.Cells.Item(5, 7).Value2 = 1.740674103860016
$ws.Cells.Item(5, 8).Value2 = 1.544363409863308
$ws.Cells.Item(5, 9).Value2 = 1.681141147197813
$ws.Cells.Item(5, 10).Value2 = 0.2179925006093075
$ws.Cells.Item(5, 11).Value2 = 0.4121521929483833
$ws.Cells.Item(5, 12).Value2 = 0.2704357263824591
$ws.Cells.Item(5, 13).Value2 = 0.1962798022966084
$ws.Cells.Item(6, 2).Value2 = 0.5753836089611326
$ws.Cells.Item(6, 4).Value2 = 0.08848707143339318
$ws.Cells.Item(6, 5).Value2 = 0.1387632501991482
$ws.Cells.Item(6, 6).Value2 = 2.435166955420542
$ws.Cells.Item(6, 7).Value2 = 1.740785286845338
$ws.Cells.Item(6, 8).Value2 = 1.54465671464763
$ws.Cells.Item(6, 9).Value2 = 1.681633374649294
$ws.Cells.Item(6, 10).Value2 = 0.2179969456063393
$ws.Cells.Item(6, 11).Value2 = 0.4101664773122593
$ws.Cells.Item(6, 12).Value2 = 0.2701756739625694
$ws.Cells.Item(6, 13).Value2 = 0.196066536047681
$ws.Cells.Item(7, 2).Value2 = 0.5801602827856414
$ws.Cells.Item(7, 4).Value2 = 0.08849610956337983
$ws.Cells.Item(7, 5).Value2 = 0.1387017015420337
$ws.Cells.Item(7, 6).Value2 = 2.433164007941059
$ws.Cells.Item(7, 7).Value2 = 1.740091096724328
$ws.Cells.Item(7, 8).Value2 = 1.542673731904699
$ws.Cells.Item(7, 9).Value2 = 1.678268286671177
$ws.Cells.Item(7, 10).Value2 = 0.2179668794109411
$ws.Cells.Item(7, 11).Value2 = 0.4239617191657317
$ws.Cells.Item(7, 12).Value2 = 0.2719923542070006
$ws.Cells.Item(7, 13).Value2 = 0.1975548400395404
$ws.Cells.Item(8, 2).Value2 = 0.6017913432446562
$ws.Cells.Item(8, 4).Value2 = 0.08856062768450634
$ws.Cells.Item(8, 5).Value2 = 0.1384513580969952
$ws.Cells.Item(8, 6).Value2 = 2.426634107721767
$ws.Cells.Item(8, 7).Value2 = 1.738946885485078
$ws.Cells.Item(8, 8).Value2 = 1.535230902027976
$ws.Cells.Item(8, 9).Value2 = 1.66469287880809
$ws.Cells.Item(8, 10).Value2 = 0.2178536297675553
$ws.Cells.Item(8, 11).Value2 = 0.4851363476130643
$ws.Cells.Item(8, 12).Value2 = 0.2802960965557872
$ws.Cells.Item(8, 13).Value2 = 0.2043196454051639
$ws.Cells.Item(9, 2).Value2 = 0.6461706485561933
$ws.Cells.Item(9, 4).Value2 = 0.08877494113198736
$ws.Cells.Item(9, 5).Value2 = 0.1380365518848408
$ws.Cells.Item(9, 6).Value2 = 2.422119327678516
$ws.Cells.Item(9, 7).Value2 = 1.74360580735906
$ws.Cells.Item(9, 8).Value2 = 1.52532951122798
$ws.Cells.Item(9, 9).Value2 = 1.642645068646146
$ws.Cells.Item(9, 10).Value2 = 0.2177009397258551
$ws.Cells.Item(9, 11).Value2 = 0.6061547364183468
$ws.Cells.Item(9, 12).Value2 = 0.2976060415829238
$ws.Cells.Item(9, 13).Value2 = 0.2182876584109259
$ws.Cells.Item(10, 2).Value2 = 0.6799539694451084
$ws.Cells.Item(10, 4).Value2 = 0.08898532985164209
$ws.Cells.Item(10, 5).Value2 = 0.1377779450797203
$ws.Cells.Item(10, 6).Value2 = 2.423853590564505
$ws.Cells.Item(10, 7).Value2 = 1.751245561633638
$ws.Cells.Item(10, 8).Value2 = 1.520912779706535
$ws.Cells.Item(10, 9).Value2 = 1.62923010860451
$ws.Cells.Item(10, 10).Value2 = 0.2176307061188485
$ws.Cells.Item(10, 11).Value2 = 0.6957041727360433
$ws.Cells.Item(10, 12).Value2 = 0.3109464729559335
$ws.Cells.Item(10, 13).Value2 = 0.2289737775409719
$ws.Cells.Item(11, 2).Value2 = 0.6955751953559002
$ws.Cells.Item(11, 4).Value2 = 0.08909241799818091
$ws.Cells.Item(11, 5).Value2 = 0.1376702679996811
$ws.Cells.Item(11, 6).Value2 = 2.425739042971031
$ws.Cells.Item(11, 7).Value2 = 1.755638829520052
$ws.Cells.Item(11, 8).Value2 = 1.51952335590093
$ws.Cells.Item(11, 9).Value2 = 1.623730792392571
$ws.Cells.Item(11, 10).Value2 = 0.2176077914787058
$ws.Cells.Item(11, 11).Value2 = 0.7365802532454495
$ws.Cells.Item(11, 12).Value2 = 0.3171501301179518
$ws.Cells.Item(11, 13).Value2 = 0.2339263789023391
$ws.Cells.Item(12, 2).Value2 = 0.7015265606635239
$ws.Cells.Item(12, 4).Value2 = 0.08913459579116889
$ws.Cells.Item(12, 5).Value2 = 0.1376309224017853
$ws.Cells.Item(12, 6).Value2 = 2.426610614251643
$ws.Cells.Item(12, 7).Value2 = 1.757434554536417
$ws.Cells.Item(12, 8).Value2 = 1.51908626698301
$ws.Cells.Item(12, 9).Value2 = 1.621735015129126
$ws.Cells.Item(12, 10).Value2 = 0.2176004074895257
$ws.Cells.Item(12, 11).Value2 = 0.7520787699467917
$ws.Cells.Item(12, 12).Value2 = 0.3195186339691958
$ws.Cells.Item(12, 13).Value2 = 0.2358148521628607
$ws.Cells.Item(13, 2).Value2 = 0.7002432359065551
$ws.Cells.Item(13, 4).Value2 = 0.08912543990872734
$ws.Cells.Item(13, 5).Value2 = 0.1376393326633233
$ws.Cells.Item(13, 6).Value2 = 2.426415898256138
$ws.Cells.Item(13, 7).Value2 = 1.757041936762221
$ws.Cells.Item(13, 8).Value2 = 1.519176442170618
$ws.Cells.Item(13, 9).Value2 = 1.622160986461509
$ws.Cells.Item(13, 10).Value2 = 0.2176019403376444
$ws.Cells.Item(13, 11).Value2 = 0.7487400209193993
$ws.Cells.Item(13, 12).Value2 = 0.3190076769514434
$ws.Cells.Item(13, 13).Value2 = 0.2354075579775383
$ws.Cells.Item(14, 2).Value2 = 0.6960640998520091
$ws.Cells.Item(14, 4).Value2 = 0.08909585546738441
$ws.Cells.Item(14, 5).Value2 = 0.1376670023865554
$ws.Cells.Item(14, 6).Value2 = 2.425807589800982
$ws.Cells.Item(14, 7).Value2 = 1.755783917070488
$ws.Cells.Item(14, 8).Value2 = 1.519485611900734
$ws.Cells.Item(14, 9).Value2 = 1.623564861379968
$ws.Cells.Item(14, 10).Value2 = 0.2176071581071621
$ws.Cells.Item(14, 11).Value2 = 0.7378549346011596
$ws.Cells.Item(14, 12).Value2 = 0.3173446018731312
$ws.Cells.Item(14, 13).Value2 = 0.234081484250531
$ws.Cells.Item(15, 2).Value2 = 0.6935089266150669
$ws.Cells.Item(15, 4).Value2 = 0.08907794555880244
$ws.Cells.Item(15, 5).Value2 = 0.1376841369354864
$ws.Cells.Item(15, 6).Value2 = 2.425455503755757
$ws.Cells.Item(15, 7).Value2 = 1.755030548531096
$ws.Cells.Item(15, 8).Value2 = 1.519686583015385
$ws.Cells.Item(15, 9).Value2 = 1.624436064081259
$ws.Cells.Item(15, 10).Value2 = 0.217610522392512
$ws.Cells.Item(15, 11).Value2 = 0.7311900485719889
$ws.Cells.Item(15, 12).Value2 = 0.3163284325022602
$ws.Cells.Item(15, 13).Value2 = 0.2332709193042106
$ws.Cells.Item(16, 2).Value2 = 0.6789381395506098
$ws.Cells.Item(16, 4).Value2 = 0.08897855948284672
$ws.Cells.Item(16, 5).Value2 = 0.1377851822391901
$ws.Cells.Item(16, 6).Value2 = 2.423752429220428
$ws.Cells.Item(16, 7).Value2 = 1.75097693006542
$ws.Cells.Item(16, 8).Value2 = 1.521016047472287
$ws.Cells.Item(16, 9).Value2 = 1.629601636048854
$ws.Cells.Item(16, 10).Value2 = 0.2176323849417408
$ws.Cells.Item(16, 11).Value2 = 0.6930355959317751
$ws.Cells.Item(16, 12).Value2 = 0.3105437576304411
$ws.Cells.Item(16, 13).Value2 = 0.2286519432095488
$ws.Cells.Item(17, 2).Value2 = 0.6700639264367112
$ws.Cells.Item(17, 4).Value2 = 0.08892049661289647
$ws.Cells.Item(17, 5).Value2 = 0.1378497198764148
$ws.Cells.Item(17, 6).Value2 = 2.422988432399535
$ws.Cells.Item(17, 7).Value2 = 1.748725332422111
$ws.Cells.Item(17, 8).Value2 = 1.521990321189932
$ws.Cells.Item(17, 9).Value2 = 1.632925011481966
$ws.Cells.Item(17, 10).Value2 = 0.2176481067256439
$ws.Cells.Item(17, 11).Value2 = 0.6696645133274615
$ws.Cells.Item(17, 12).Value2 = 0.3070295571357775
$ws.Cells.Item(17, 13).Value2 = 0.2258416853630578
$ws.Cells.Item(18, 2).Value2 = 0.6649835556326025
$ws.Cells.Item(18, 4).Value2 = 0.08888817227431289
$ws.Cells.Item(18, 5).Value2 = 0.1378877783248549
$ws.Cells.Item(18, 6).Value2 = 2.422652217198106
$ws.Cells.Item(18, 7).Value2 = 1.747516664844781
$ws.Cells.Item(18, 8).Value2 = 1.522609037200255
$ws.Cells.Item(18, 9).Value2 = 1.634893309500605
$ws.Cells.Item(18, 10).Value2 = 0.2176580002431052
$ws.Cells.Item(18, 11).Value2 = 0.6562352500783106
$ws.Cells.Item(18, 12).Value2 = 0.3050209975707503
$ws.Cells.Item(18, 13).Value2 = 0.2242339123776134
$ws.Cells.Item(19, 2).Value2 = 0.6632675376305883
$ws.Cells.Item(19, 4).Value2 = 0.08887741223842482
$ws.Cells.Item(19, 5).Value2 = 0.1379008254900256
$ws.Cells.Item(19, 6).Value2 = 2.422556110577233
$ws.Cells.Item(19, 7).Value2 = 1.747122265728294
$ws.Cells.Item(19, 8).Value2 = 1.522828545725389
$ws.Cells.Item(19, 9).Value2 = 1.635569494928383
$ws.Cells.Item(19, 10).Value2 = 0.2176614963197139
$ws.Cells.Item(19, 11).Value2 = 0.6516906069333857
$ws.Cells.Item(19, 12).Value2 = 0.3043431208924261
$ws.Cells.Item(19, 13).Value2 = 0.223691030942561
$ws.Cells.Item(20, 2).Value2 = 0.6710061372120038
$ws.Cells.Item(20, 4).Value2 = 0.08892656664732357
$ws.Cells.Item(20, 5).Value2 = 0.1378427526702248
$ws.Cells.Item(20, 6).Value2 = 2.42305907937579
$ws.Cells.Item(20, 7).Value2 = 1.748956076903994
$ws.Cells.Item(20, 8).Value2 = 1.521880570645308
$ws.Cells.Item(20, 9).Value2 = 1.632565356279713
$ws.Cells.Item(20, 10).Value2 = 0.2176463451074895
$ws.Cells.Item(20, 11).Value2 = 0.6721510456127078
$ws.Cells.Item(20, 12).Value2 = 0.3074023345100727
$ws.Cells.Item(20, 13).Value2 = 0.2261399515463296
$ws.Cells.Item(21, 2).Value2 = 0.6972906413317901
$ws.Cells.Item(21, 4).Value2 = 0.08910450109461721
$ws.Cells.Item(21, 5).Value2 = 0.1376588363528342
$ws.Cells.Item(21, 6).Value2 = 2.425981988452918
$ws.Cells.Item(21, 7).Value2 = 1.756149842302818
$ws.Cells.Item(21, 8).Value2 = 1.519392384834248
$ws.Cells.Item(21, 9).Value2 = 1.623150156798623
$ws.Cells.Item(21, 10).Value2 = 0.2176055904691943
$ws.Cells.Item(21, 11).Value2 = 0.7410516190020928
$ws.Cells.Item(21, 12).Value2 = 0.3178325639937896
$ws.Cells.Item(21, 13).Value2 = 0.234470631532993
$ws.Cells.Item(22, 2).Value2 = 0.71467840359918
$ws.Cells.Item(22, 4).Value2 = 0.08923026333528128
$ws.Cells.Item(22, 5).Value2 = 0.1375469663395443
$ws.Cells.Item(22, 6).Value2 = 2.428810801266906
$ws.Cells.Item(22, 7).Value2 = 1.761621349028673
$ws.Cells.Item(22, 8).Value2 = 1.518285266314052
$ws.Cells.Item(22, 9).Value2 = 1.617502079575864
$ws.Cells.Item(22, 10).Value2 = 0.2175864905592597
$ws.Cells.Item(22, 11).Value2 = 0.7861963551794418
$ws.Cells.Item(22, 12).Value2 = 0.3247618484069079
$ws.Cells.Item(22, 13).Value2 = 0.2399911108685231
$ws.Cells.Item(23, 2).Value2 = 0.7053792191911441
$ws.Cells.Item(23, 4).Value2 = 0.08916227848709113
$ws.Cells.Item(23, 5).Value2 = 0.1376059124148543
$ws.Cells.Item(23, 6).Value2 = 2.427216994640631
$ws.Cells.Item(23, 7).Value2 = 1.758630619277653
$ws.Cells.Item(23, 8).Value2 = 1.518828684465518
$ws.Cells.Item(23, 9).Value2 = 1.620470345225542
$ws.Cells.Item(23, 10).Value2 = 0.2175959969618493
$ws.Cells.Item(23, 11).Value2 = 0.7620914529508696
$ws.Cells.Item(23, 12).Value2 = 0.3210532966979684
$ws.Cells.Item(23, 13).Value2 = 0.2370378210249626
$ws.Cells.Item(24, 2).Value2 = 0.6705800966165327
$ws.Cells.Item(24, 4).Value2 = 0.08892381909284097
$ws.Cells.Item(24, 5).Value2 = 0.1378458995701394
$ws.Cells.Item(24, 6).Value2 = 2.42302681897327
$ws.Cells.Item(24, 7).Value2 = 1.748851490045695
$ws.Cells.Item(24, 8).Value2 = 1.52193000634054
$ws.Cells.Item(24, 9).Value2 = 1.632727776826798
$ws.Cells.Item(24, 10).Value2 = 0.2176471388725556
$ws.Cells.Item(24, 11).Value2 = 0.6710268622722708
$ws.Cells.Item(24, 12).Value2 = 0.3072337650922634
$ws.Cells.Item(24, 13).Value2 = 0.22600508084561
$ws.Cells.Item(25, 2).Value2 = 0.6339567698218218
$ws.Cells.Item(25, 4).Value2 = 0.08870762653446107
$ws.Cells.Item(25, 5).Value2 = 0.1381406454595373
$ws.Cells.Item(25, 6).Value2 = 2.422453443993533
$ws.Cells.Item(25, 7).Value2 = 1.741605557161023
$ws.Cells.Item(25, 8).Value2 = 1.527505969003428
$ws.Cells.Item(25, 9).Value2 = 1.648120433597718
$ws.Cells.Item(25, 10).Value2 = 0.2177348547038669
$ws.Cells.Item(25, 11).Value2 = 0.5733037300140325
$ws.Cells.Item(25, 12).Value2 = 0.2928136657487954
$ws.Cells.Item(25, 13).Value2 = 0.2144342059729638
